$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.955.05"
$ws.Range("E2").Value = "  +5.97%  "

$ws.Range("D3").Value = "3.660.21"
$ws.Range("E3").Value = "  +17.98%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.85"
$ws.Range("E5").Value = "  +7.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.68"
$ws.Range("E6").Value = "  +1.98%  "

$ws.Range("D7").Value = "3.657.45"
$ws.Range("E7").Value = "  +17.93%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +4.81%  "

$ws.Range("E10").Value = "  +8.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("E11").Value = "  +5.67%  "

$ws.Range("E12").Value = "  +7.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.46"
$ws.Range("E13").Value = "  +11.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000254"
$ws.Range("E14").Value = "  +5.80%  "

$ws.Range("D15").Value = "4.268.38"
$ws.Range("E15").Value = "  +17.96%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.674.86"
$ws.Range("E16").Value = "  +18.55%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "70.949.97"
$ws.Range("E17").Value = "  +6.08%  "

$ws.Range("E18").Value = "  +1.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("E19").Value = "  +7.12%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "518.90"
$ws.Range("E20").Value = "  +7.86%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.87"
$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("E22").Value = "  +19.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.743"
$ws.Range("E23").Value = "  +7.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.72"
$ws.Range("E24").Value = "  +5.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  +11.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.51"
$ws.Range("E26").Value = "  +6.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.08"
$ws.Range("E27").Value = "  +9.37%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("E29").Value = "  +12.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("E30").Value = "  +12.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.09"
$ws.Range("E31").Value = "  +1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.70"
$ws.Range("E32").Value = "  +13.14%  "

$ws.Range("E33").Value = "  +17.16%  "

$ws.Range("E34").Value = "  +3.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.13"
$ws.Range("E36").Value = "  +9.62%  "

$ws.Range("E37").Value = "  +8.64%  "

$ws.Range("E38").Value = "  +11.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  +10.33%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.132"
$ws.Range("E40").Value = "  +6.86%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.81"
$ws.Range("E41").Value = "  +5.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "46.03"
$ws.Range("E42").Value = "  -4.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  +5.97%  "

$ws.Range("D44").Value = "3.120.57"
$ws.Range("E44").Value = "  +11.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "423.57"
$ws.Range("E45").Value = "  +13.24%  "

$ws.Range("E46").Value = "  +3.99%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0371"
$ws.Range("E47").Value = "  +8.21%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.51"
$ws.Range("E48").Value = "  +12.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.59"
$ws.Range("E49").Value = "  +3.90%  "

$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("E51").Value = "  +10.05%  "
